$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10-24 down to 11-25.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly price record.
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44601
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 100112045
$ws.Range("G10").Value = "Zapallo"
$ws.Range("H10").Value = "Camote"
$ws.Range("I10").Value = "2a (cosecha)"
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 400
$ws.Range("L10").Value = 450
$ws.Range("M10").Value = 425
$ws.Range("N10").Value = "$/kilo (volumen en unidades)"
$ws.Range("O10").Value = "Región de O'Higgins"
$ws.Range("P10").Value = 425
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"
